$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I9").Value = 2
$ws.Range("G10").Value = 2
$ws.Range("K11").Value = 2
$ws.Range("F14").Value = 1

$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 1
$ws.Range("J16").Value = 0

$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 2
$ws.Range("G17").Value = 1
$ws.Range("H17").Value = 1
$ws.Range("I17").Value = 1

$ws.Range("O17").Select()
